# This workbook (Location.xlsx) was opened in Excel and re-saved without any
# underlying data edits (commit "Multimedia and Spreadsheet Data, script
# files" touches a batch of repo assets). The diff is Excel's own
# normalization of the OOXML package: workbook/sheet view bookkeeping,
# style-table canonicalisation, explicit row heights, and dropping the
# (already content-less) drawing part. Reproduce the observable, user-facing
# pieces of that resave through the object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View state: 150% zoom, active cell/selection parked at E36 ----------
$excel.ActiveWindow.Zoom = 150
$null = $ws.Range("E36").Select()

# --- Column widths (B:D) widened slightly, as in the target sheet --------
# (ColumnWidth uses Excel's character-unit scale; the inputs below are the
# values that round-trip to the target XML `width` attributes of 25.5,
# 24.83203125 and 16 for columns B, C and D respectively.)
$ws.Columns.Item(2).ColumnWidth = 24.609375
$ws.Columns.Item(3).ColumnWidth = 23.916666666666664
$ws.Columns.Item(4).ColumnWidth = 15.234375

# --- Row heights: every populated row (1-26) gets an explicit 15.75pt ----
# height, matching the `ht="15.75" customHeight="1"` Excel now stamps on
# each <row> on save.
for ($r = 1; $r -le 26; $r++) {
  $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Drop the (empty) drawing object ---------------------------------
# The original workbook carries a <xdr:wsDr/> drawing part with zero
# shapes/charts/pictures in it; Excel's resave drops that now-pointless
# part entirely. Make sure there truly is nothing left referencing it.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
  $ws.Shapes.Item($i).Delete()
}
